$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to add for columns S (2021) and T (2022), keyed by row number.
$values = @{
    4  = @(2021, 2022)
    5  = @(2.5, 2.6)
    6  = @(2.5, 1.8)
    7  = @(1.6, 2.6)
    8  = @(3.6, 1.9)
    9  = @(5.8, 3.9)
    10 = @(1.1000000000000001, 3.2)
    11 = @(1.1000000000000001, 3.3)
    12 = @(5.0999999999999996, 2.5)
    13 = @(2.2999999999999998, 1.9)
    14 = @(2.1, 2.5)
}

foreach ($row in 4..14) {
    $r = $ws.Range("R$row")

    # Copy the formatting (cell style) of column R into the new S/T cells,
    # then overwrite the pasted values with the actual data for that column.
    $r.Copy() | Out-Null
    $s = $ws.Range("S$row")
    $s.PasteSpecial(-4122) | Out-Null
    $t = $ws.Range("T$row")
    $t.PasteSpecial(-4122) | Out-Null

    $vals = $values[$row]
    $s.Value = $vals[0]
    $t.Value = $vals[1]
}

# Update the selected cell in the sheet view, matching the new selection.
$ws.Range("V7").Select() | Out-Null
